$d = $word.ActiveDocument

# Remove " and Python 2.7" from the "Primary language: PHP 6.0 ..." bullet,
# while leaving the existing _GoBack bookmark (which sits between
# "Python" and " 2.7") untouched. Doing this as two separate
# Find/Replace passes that each stop short of the bookmark keeps it intact.
$d.Content.Find.Execute(" and Python", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
$d.Content.Find.Execute(" 2.7", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# Add a new bulleted line "Framework: Yii 1.1.16" right after the
# "IDE: PhpStorm 10" bullet, inheriting the same list paragraph formatting.
$d.Content.Find.Execute("IDE: PhpStorm 10", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "IDE: PhpStorm 10^pFramework: Yii 1.1.16", 2)
